# Remove the trailing "Source: <hyperlink>" paragraph (and the now-
# superfluous blank paragraph that immediately preceded it) from the
# end of the document, as part of refreshing the assets for the next
# academic year's course iteration.

$d = $word.ActiveDocument

$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)

if ($lastPara.Range.Text -like "Source:*") {
    # Delete the "Source: https://..." paragraph entirely (text, run
    # properties and the hyperlink field all go with it).
    $lastPara.Range.Delete()

    # The paragraph that used to sit right above it was just an empty
    # spacer paragraph; remove it too so the trailing blank-paragraph
    # count goes back down to what it was before that source line was
    # appended.
    $newCount = $d.Paragraphs.Count
    $spacerPara = $d.Paragraphs.Item($newCount - 1)
    if ($spacerPara.Range.Text.Trim() -eq "") {
        $spacerPara.Range.Delete()
    }
}
